$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Column B (accuracy values) per diff
$ws.Range("B3").Value = 0.203125
$ws.Range("B4").Value = 0.125
$ws.Range("B5").Value = 0.09375
$ws.Range("B6").Value = 0.125
$ws.Range("B7").Value = 0.109375
$ws.Range("B8").Value = 0.109375
$ws.Range("B9").Value = 0.109375
$ws.Range("B10").Value = 0.09375
$ws.Range("B11").Value = 0.09375
$ws.Range("B12").Value = 0.109375
$ws.Range("B13").Value = 0.09375
$ws.Range("B14").Value = 0.09375
$ws.Range("B15").Value = 0.09375
$ws.Range("B17").Value = 0.109375
$ws.Range("B18").Value = 0.109375
$ws.Range("B19").Value = 0.109375
$ws.Range("B20").Value = 0.109375
$ws.Range("B21").Value = 0.109375
$ws.Range("B22").Value = 0.109375
$ws.Range("B23").Value = 0.109375
$ws.Range("B24").Value = 0.109375
$ws.Range("B25").Value = 0.109375
$ws.Range("B26").Value = 0.109375
$ws.Range("B27").Value = 0.109375
$ws.Range("B28").Value = 0.109375
$ws.Range("B29").Value = 0.109375
$ws.Range("B30").Value = 0.109375
$ws.Range("B31").Value = 0.109375
$ws.Range("B32").Value = 0.109375
$ws.Range("B33").Value = 0.109375
$ws.Range("B34").Value = 0.09375
$ws.Range("B35").Value = 0.09375
$ws.Range("B36").Value = 0.09375
$ws.Range("B37").Value = 0.09375
$ws.Range("B38").Value = 0.09375
$ws.Range("B39").Value = 0.09375
$ws.Range("B40").Value = 0.09375
$ws.Range("B41").Value = 0.09375
$ws.Range("B42").Value = 0.09375
$ws.Range("B43").Value = 0.09375
$ws.Range("B44").Value = 0.09375
$ws.Range("B45").Value = 0.09375
$ws.Range("B46").Value = 0.09375
$ws.Range("B69").Value = 0.078125
$ws.Range("B70").Value = 0.078125
$ws.Range("B71").Value = 0.078125
$ws.Range("B72").Value = 0.078125
$ws.Range("B73").Value = 0.078125
$ws.Range("B74").Value = 0.078125
$ws.Range("B75").Value = 0.078125
$ws.Range("B76").Value = 0.078125
$ws.Range("B77").Value = 0.078125
$ws.Range("B78").Value = 0.078125
$ws.Range("B79").Value = 0.078125
$ws.Range("B80").Value = 0.078125
$ws.Range("B81").Value = 0.078125
$ws.Range("B82").Value = 0.078125
$ws.Range("B83").Value = 0.078125
$ws.Range("B84").Value = 0.078125
$ws.Range("B85").Value = 0.078125
$ws.Range("B86").Value = 0.078125
$ws.Range("B87").Value = 0.078125
$ws.Range("B88").Value = 0.0625
$ws.Range("B89").Value = 0.0625
$ws.Range("B90").Value = 0.0625
$ws.Range("B91").Value = 0.0625
$ws.Range("B92").Value = 0.0625
$ws.Range("B93").Value = 0.0625
$ws.Range("B94").Value = 0.0625
$ws.Range("B95").Value = 0.0625
$ws.Range("B96").Value = 0.0625
$ws.Range("B97").Value = 0.0625
$ws.Range("B98").Value = 0.0625
$ws.Range("B99").Value = 0.0625
$ws.Range("B100").Value = 0.0625
$ws.Range("B101").Value = 0.0625
$ws.Range("B102").Value = 0.0625
$ws.Range("B103").Value = 0.09375
$ws.Range("B104").Value = 0.140625
$ws.Range("B106").Value = 0.15625
$ws.Range("B107").Value = 0.109375
$ws.Range("B108").Value = 0.140625
$ws.Range("B109").Value = 0.109375
$ws.Range("B110").Value = 0.109375
$ws.Range("B111").Value = 0.09375
$ws.Range("B113").Value = 0.140625
$ws.Range("B114").Value = 0.09375
$ws.Range("B115").Value = 0.09375
$ws.Range("B116").Value = 0.09375
$ws.Range("B117").Value = 0.09375
$ws.Range("B118").Value = 0.180327868852459

# Update Column A inline string object addresses (rows 102-118)
$ws.Range("A102").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A103").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A104").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A105").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A106").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A107").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A108").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A109").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A110").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A111").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A112").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A113").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A114").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A115").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A116").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A117").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
$ws.Range("A118").Value = "<__main__.DisplayOutputs object at 0x7fdc0410bb20>"
